# Manna Drain mostly Done
# Update the "Data" sheet that drives the Gantt chart: several tasks have
# shifted start dates and/or flipped from "Remaining" to "Completed".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# Row 3 : Start Date 40921 -> 40925, Completed 0 -> 1, Remaining 1 -> 0
$ws.Range("B3").Value = 40925
$ws.Range("C3").Value = 1
$ws.Range("D3").Value = 0

# Row 10 : Start Date 40921 -> 40926, Completed 0 -> 1, Remaining 8 -> 0
$ws.Range("B10").Value = 40926
$ws.Range("C10").Value = 1
$ws.Range("D10").Value = 0

# Rows 15-20 : Start Date shifts only (Completed/Remaining unchanged)
$ws.Range("B15").Value = 40926
$ws.Range("B16").Value = 40926
$ws.Range("B17").Value = 40926
$ws.Range("B18").Value = 40927
$ws.Range("B19").Value = 40927
$ws.Range("B20").Value = 40927

# Row 22 : Start Date 40914 -> 40927, Completed 5 -> 0, Remaining 2 -> 1
$ws.Range("B22").Value = 40927
$ws.Range("C22").Value = 0
$ws.Range("D22").Value = 1

# Rows 24, 26, 30, 31 : Start Date shifts only
$ws.Range("B24").Value = 40928
$ws.Range("B26").Value = 40928
$ws.Range("B30").Value = 40928
$ws.Range("B31").Value = 40928

# Restore the view: scrolled up a bit and selection moved to D10
$ws.Activate() | Out-Null
$ws.Range("D10").Select() | Out-Null
